$wb = $excel.ActiveWorkbook

# --- ProviderSearch sheet: insert a new "COUNTY" column between CITY and ZIP_CODE ---
$ws = $wb.Worksheets.Item("ProviderSearch")

# Inserting column T shifts CITY's old neighbour (ZIP_CODE, originally column T) and
# everything to its right one column to the right (T->U, U->V, ... AG->AH), carrying
# values/styles/column-widths along - exactly like Excel's native "Insert Column".
$ws.Columns("T").Insert() | Out-Null

# New header + data for the inserted COUNTY column.
$ws.Range("T1").Value = "COUNTY"
$ws.Range("T2:T11").Value = "Alpine"

# Make ProviderSearch the active/selected sheet and tab, with the same selected cell
# that the authored workbook ends up with.
$ws.Activate() | Out-Null
$ws.Range("T15").Select() | Out-Null
